# Pol Mandated Cap Const Sched.xlsx - apply Jun's updates:
#  - About sheet: replace the first "Source" block (Virginia Clean Economy Act
#    citation + second Dominion/Greentech source block) with a single
#    "None (this variable is intended to be user-specified)" note, removing
#    the now-unused rows.
#  - Delete the "Required offshore wind" sheet entirely.
#  - PMCCS sheet: the array formula in B14 referenced the now-deleted sheet,
#    so clear it (values fall back to 0, matching what Excel leaves behind
#    once the precedent sheet disappears).

$wb = $excel.ActiveWorkbook

# --- About sheet -----------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Replace the "Source:" value in B3 with the new placeholder note. It
# loses the bold/shaded "source header" look that went with the old
# two-column source layout, and the shaded C3:D3 filler cells go away too.
$about.Range("B3").Value = "None (this variable is intended to be user-specified)"
$about.Range("B3").ClearFormats()
$about.Range("C3:D3").Clear()

# Remove the rest of the old source block (rows 4-12): second source
# citation, dates, links, etc. Everything below (the "Note:" block) shifts
# up to take their place.
$about.Rows("4:12").Delete()

# --- PMCCS sheet -------------------------------------------------------
$pmccs = $wb.Worksheets.Item("PMCCS")

# This array formula pulled the construction schedule from the
# "Required offshore wind" sheet; replace it with the literal zeros that
# Excel leaves behind once that sheet (and the formula with it) is gone.
$pmccs.Range("B14:AI14").Value = 0

# --- Remove the "Required offshore wind" sheet --------------------------
$wb.Worksheets.Item("Required offshore wind").Delete() | Out-Null
